$d = $word.ActiveDocument

# Locate the closing paragraph of "2. Problemas del sistema Ecobici", which
# originally begins "Desde la primera etapa del sistema..." and ends
# "...de mayor demanda de candados. " (right before the _GoBack bookmark).
$find = $d.Content.Find
$found = $find.Execute(
    "Desde la primera etapa del sistema*de mayor demanda de candados. ",
    $false, $false, $true, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate target paragraph text"
}

$r = $find.Parent

# Rewrite the whole paragraph body in place. Assigning to Range.Text keeps
# the run formatting (Times New Roman, color, size, es-ES) that the source
# range already carried, and leaves the trailing bookmark (_GoBack) and the
# paragraph mark that follow untouched.
$r.Text = "La primera etapa del sistema ayudó a identificar el principal problema que existe, este es la falta de bicicletas o espacios para estacionar las unidades en algunas de las estaciones; se presenta porque las estaciones, durante ciertas horas en el día, tienen distinta demanda y se comportan de distinta forma. De acuerdo a su demanda, se pueden considerar tres tipos de estaciones: demandan bicicletas, son las estaciones en las que hay demasiados espacios vacíos y pocas bicicletas durante un periodo largo de tiempo; demandan espacios, estaciones que durante un periodo de tiempo tienen demasiadas bicicletas y pocos espacios; autobalanceables "
